$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix spelling error: "Offsett" -> "Offset"
$ws.Range("G1").Value = "Offset"

# Update the F2 value to include the "-[id]" suffix
$ws.Range("F2").Value = "CaliSimulationResults/water-level-reach-0-[id]"

# Update the active selection to F2 (matches topLeftCell/selection change)
$ws.Range("F2").Select()

$wb.Save()
